$d = $word.ActiveDocument

# Remove the three trailing paragraphs that follow the
# "LOQ4083: Fenômenos de Transporte I (Requisito fraco)" paragraph:
#   1) an empty paragraph
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github
#       pages. Original theme under Creative Commons Attribution"
# The empty paragraph and the page-break paragraph that follow must stay.

$target = $d.Content.Find
$target.ClearFormatting()
$target.Text = "Ver no Jupiter Salvar em pdf Salvar em docx"
$target.Execute() | Out-Null

if ($target.Found) {
    $r = $target.Parent
    # Extend the range to cover the preceding empty paragraph, this
    # paragraph, and the following "© 2020 ..." paragraph (3 paragraphs
    # total), then delete them all in one shot.
    $startPara = $r.Paragraphs(1)
    $prevPara = $startPara.Previous()
    $nextPara = $startPara.Next()

    $delStart = $prevPara.Range.Start
    $delEnd = $nextPara.Range.End

    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
